# Update the "District" column (G) with the official district names,
# and drop the two stray empty "Address" cells (F28/F33) that the
# original sheet carried as blank inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected district name (column G)
$districtUpdates = @{
    3  = "Vijayapura (Bijapur)"
    4  = "Uttara Kannada (Karwar)"
    5  = "Vijayapura (Bijapur)"
    6  = "Vijayapura (Bijapur)"
    7  = "Vijayapura (Bijapur)"
    9  = "Uttara Kannada (Karwar)"
    10 = "Vijayapura (Bijapur)"
    11 = "Uttara Kannada (Karwar)"
    12 = "Uttara Kannada (Karwar)"
    13 = "Uttara Kannada (Karwar)"
    15 = "Vijayapura (Bijapur)"
    16 = "Uttara Kannada (Karwar)"
    17 = "Vijayapura (Bijapur)"
    18 = "Uttara Kannada (Karwar)"
    19 = "Vijayapura (Bijapur)"
    20 = "Uttara Kannada (Karwar)"
    22 = "Vijayapura (Bijapur)"
    26 = "Uttara Kannada (Karwar)"
    27 = "Vijayapura (Bijapur)"
    29 = "Vijayapura (Bijapur)"
    30 = "Uttara Kannada (Karwar)"
    31 = "Uttara Kannada (Karwar)"
    32 = "Vijayapura (Bijapur)"
    34 = "Uttara Kannada (Karwar)"
    35 = "Vijayapura (Bijapur)"
}

foreach ($row in $districtUpdates.Keys) {
    $ws.Range("G$row").Value = $districtUpdates[$row]
}

# Rows 28 and 33 have a leftover empty "Address" cell (F) with no text;
# clear it so the cell is dropped entirely, matching the official list.
$ws.Range("F28").ClearContents()
$ws.Range("F33").ClearContents()
